$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 4
$ws.Range("D1").Value = 40
$ws.Range("F1").Value = 40
$ws.Range("H1").Value = 4

$ws.Range("F3").Select()
$ws.Range("F3:F91").Select()
